$wb = $excel.ActiveWorkbook

$wsCompounds = $wb.Worksheets.Item("compounds")
$wsBiomarkers = $wb.Worksheets.Item("biomarkers")

# Update version values on "compounds" sheet
# Row 2 = Open Targets Platform -> version bump 2023.02 -> 2023.06
$wsCompounds.Range("E2").NumberFormat = "@"
$wsCompounds.Range("E2").Value = "2023.06"
# Row 3 = NCI Thesaurus -> version bump 23.05e -> 23.06d
$wsCompounds.Range("E3").Value = "23.06d"

# Update selection / active cell on compounds sheet
$wsCompounds.Range("E2").Select()

# Update selection / active cell on biomarkers sheet (A4:XFD4 stays, topLeftCell changes)
$wsBiomarkers.Activate()
$excel.ActiveWindow.ScrollColumn = 2

$wsBiomarkers.Range("A4:XFD4").Select()
